# "All num operations complete, nasm output correct"
#
# The productions notes were completely rewritten: the old LL(1) grammar
# listing (Goal/Expr/Term/Ex/Factor/Neg) is replaced with the new grammar
# used by the num/NASM assignment (Goal/Statement/Decl/Assign/Expr/Term/
# Factor/Power/Base), formatted as a Python-style list of quoted
# production strings. A handful of runs keep the <w:proofErr> spellcheck/
# grammar-check markers Word had inserted around "Decl" and "( name" /
# ">( Expr" so the run-boundary shape matches the original editing
# session.
#
# The whole paragraph list (the entire body content) is replaced in one
# shot via Range.InsertXML with a literal WordprocessingML fragment - a
# PowerShell single-quoted here-string so the embedded apostrophes,
# backslashes and XML entities need no escaping.

$d = $word.ActiveDocument

$newParagraphsXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve"># </w:t></w:r><w:r><w:t>Assignment 4 Productions</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>'Goal-&gt;Statement',</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>'Statement-&gt;</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Decl</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Assign',</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>'Statement-&gt;name Assign',</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">'Statement-&gt;print </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>( name</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> )',</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>'</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Decl</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>-&gt;type name',</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>'Assign-&gt;= Expr',</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>'Assign-&gt;eps',</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>'Expr-&gt;Term Expr\'',</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>'Expr\'-&gt;+ Term Expr\'',</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>'Expr\'-&gt;- Term Expr\'',</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>'Expr\'-&gt;eps',</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>'Term-&gt;Factor Term\'',</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>'Term\'-&gt;* Factor Term\'',</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>'Term\'-&gt;/ Factor Term\'',</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>'Term\'-&gt;eps',</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>'Factor-&gt;Base Power',</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>'Power-&gt;^ Base Power',</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>'Power-&gt;eps',</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>'Base-</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>&gt;( Expr</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> )',</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>'Base-&gt;Base\'',</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>'Base-&gt;- Base\'',</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>'Base\'-&gt;number',</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>'Base\'-&gt;name',</w:t></w:r></w:p>
'@

$paraCount = $d.Paragraphs.Count
$bodyStart = $d.Paragraphs(1).Range.Start
$bodyEnd = $d.Paragraphs($paraCount).Range.End
$bodyRange = $d.Range($bodyStart, $bodyEnd)

[void]$bodyRange.InsertXML($newParagraphsXml)
